# Update the cryptocurrency price/volume snapshot (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> cell edits. "D" (Price) values must be written as literal text so
# that things like trailing zeros ("25.20"), thousand-dot grouping
# ("28.690.12") and subscript-digit notation ("0.0" + subscript-3 + "0690")
# survive exactly as authored instead of being auto-coerced into numbers.
$updates = @(
    @{ Row = 2;  D = "28.690.12";  E = "  +1.46%  " }
    @{ Row = 3;  D = "1.571.43";   E = "  +1.28%  " }
    @{ Row = 4;  D = "0.997";      E = "  -0.58%  " }
    @{ Row = 5;  D = "210.74";     E = "  +0.78%  " }
    @{ Row = 6;  D = "0.517";      E = "  +6.59%  " }
    @{ Row = 7;  D = "0.996";      E = "  -0.68%  " }
    @{ Row = 8;  D = "25.20";      E = "  +7.57%  " }
    @{ Row = 9;                    E = "  +1.89%  " }
    @{ Row = 10; D = "0.0589";     E = "  +1.23%  " }
    @{ Row = 11; D = "0.0900";     E = "  +1.15%  " }
    @{ Row = 12; D = "1.790.64";   E = "  +0.94%  " }
    @{ Row = 13; D = "1.567.14";   E = "  +1.00%  " }
    @{ Row = 14; D = "28.708.56";  E = "  +1.49%  " }
    @{ Row = 15; D = "0.519";      E = "  +2.23%  " }
    @{ Row = 16; D = "3.67";       E = "  +1.20%  " }
    @{ Row = 17; D = "61.87";      E = "  +2.68%  " }
    @{ Row = 18; D = "229.88";     E = "  +0.89%  " }
    @{ Row = 19; D = "7.32";       E = "  +0.45%  " }
    @{ Row = 20; D = "0.0$([char]0x2083)0690"; E = "  +2.80%  " }
    @{ Row = 21; D = "0.995";      E = "  -0.78%  " }
    @{ Row = 22; D = "3.96";       E = "  +1.41%  " }
    @{ Row = 23; D = "9.12";       E = "  +3.62%  " }
    @{ Row = 24; D = "2.08";       E = "  +3.69%  " }
    @{ Row = 25; D = "152.46";     E = "  +3.43%  " }
    @{ Row = 26;                   E = "  +3.96%  " }
    @{ Row = 27; D = "14.88";      E = "  +1.00%  " }
    @{ Row = 28; D = "6.28";       E = "  +1.10%  " }
    @{ Row = 29; D = "0.997";      E = "  -0.64%  " }
    @{ Row = 30; D = "0.0461";     E = "  -1.19%  " }
    @{ Row = 31; D = "1.06";       E = "  -0.47%  " }
    @{ Row = 32; D = "3.20";       E = "  +1.16%  " }
    @{ Row = 33; D = "1.410.37";   E = "  +1.97%  " }
    @{ Row = 34; D = "3.02";       E = "  -0.58%  " }
    @{ Row = 35;                   E = "  -2.10%  " }
    @{ Row = 36; D = "1.49";       E = "  -0.21%  " }
    @{ Row = 37;                   E = "  +5.84%  " }
    @{ Row = 38;                   E = "  -1.64%  " }
    @{ Row = 39;                   E = "  +0.89%  " }
    @{ Row = 40; B = "ImmutableX";  C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx";          D = "0.520"; E = "  +2.05%  " }
    @{ Row = 41; B = "RenderToken"; C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr";    D = "1.95";  E = "  +1.73%  " }
    @{ Row = 42; D = "0.997";      E = "  -0.61%  " }
    @{ Row = 43; D = "0.776";      E = "  +0.50%  " }
    @{ Row = 44; D = "0.0462";     E = "  -0.46%  " }
    @{ Row = 45; D = "63.64";      E = "  +3.37%  " }
    @{ Row = 46; D = "5.27";       E = "  -2.06%  " }
    @{ Row = 47; D = "1.705.76";   E = "  +1.12%  " }
    @{ Row = 48; B = "Quant";      C = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt";            D = "84.82"; E = "  -0.32%  " }
    @{ Row = 49; B = "WEMIXToken"; C = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix";          D = "0.825"; E = "  -9.32%  " }
    @{ Row = 50; B = "BitcoinSV";  C = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv";         D = "42.14"; E = "  +1.67%  " }
    @{ Row = 51; B = "BabyDogeCoin"; C = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge";  D = "0.0$([char]0x2086)0102"; E = "  -1.39%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey("B")) {
        $ws.Range("B$row").Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $ws.Range("C$row").Value = $u.C
    }
    if ($u.ContainsKey("D")) {
        # Force text so numeric-looking prices keep their exact original
        # representation (trailing zeros, dotted thousands, subscripts...).
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E$row").Value = $u.E
    }
}
